$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999764011405,
    0.99824553993083698,
    0.994180289614639,
    0.99704366128309929,
    0.98675255465820111,
    0.96298321046663449,
    0.95812779374344426,
    0.95074632253765012,
    0.9424140113314885,
    0.93499040034144354,
    0.93386517924059231,
    0.93198047171457887,
    0.92567992889903628,
    0.92397995554556689,
    0.92351519942926164,
    0.92342650054425246,
    0.92561994839808404,
    0.92451101863433061,
    0.99695292417857562,
    0.98983574454029077,
    0.98843722668880452,
    0.9871727176618168,
    0.9871683994531999,
    0.97421412483617886,
    0.9678216128834467,
    0.96081298733393083,
    0.95803652196062428,
    0.94720656711446627,
    0.93983046396043168,
    0.9341944423017442,
    0.93006655579481567,
    0.92971558375354713,
    0.93115203010566727
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

